$d = $word.ActiveDocument

$d.Content.Find.Execute("Subject line:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "الموضوع:", 2)

$d.Content.Find.Execute("Body:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "المحتوى:", 2)

$d.Content.Find.Execute("Important", $true, $false, $false, $false, $false,
                         $true, 1, $false, "هام", 2)

$d.Content.Find.Execute("Live chat", $true, $false, $false, $false, $false,
                         $true, 1, $false, "دردشة حية", 2)

$d.Content.Find.Execute("WhatsApp", $true, $false, $false, $false, $false,
                         $true, 1, $false, "واتساب", 2)
